# g8.3 - reestruturação para evitar referência
# Unpivot the table from wide format (Produto | 2025-2024 | 2025/1997)
# into a long format (Produto | Categoria | Valor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 10

# --- 1. Capture the original data (rows 2..10, columns A/B/C) before overwriting ---
$products = @()
$valsB = @()
$valsC = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $products += ,$ws.Cells.Item($r, 1).Value2
    $valsB += ,$ws.Cells.Item($r, 2).Value2
    $valsC += ,$ws.Cells.Item($r, 3).Value2
}

$count = $products.Count

# --- 2. Update header row ---
$ws.Cells.Item(1, 2).Value2 = "Categoria"
$ws.Cells.Item(1, 3).Value2 = "Valor"

# --- 3. Write block 1: category "2025-2024" using the original column B values ---
for ($i = 0; $i -lt $count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value2 = $products[$i]
    $ws.Cells.Item($r, 2).Value2 = "2025-2024"
    $v = $valsB[$i]
    if ($v -eq $null -or $v -eq "") {
        $ws.Cells.Item($r, 3).Value2 = ""
    } else {
        $ws.Cells.Item($r, 3).Value2 = $v
    }
}

# --- 4. Write block 2: category "2025/1997" using the original column C values ---
for ($i = 0; $i -lt $count; $i++) {
    $r = 2 + $count + $i
    $ws.Cells.Item($r, 1).Value2 = $products[$i]
    $ws.Cells.Item($r, 2).Value2 = "2025/1997"
    $v = $valsC[$i]
    if ($v -eq $null -or $v -eq "") {
        $ws.Cells.Item($r, 3).Value2 = ""
    } else {
        $ws.Cells.Item($r, 3).Value2 = $v
    }
}
